# Update countries & provincias Spain
# Applies the 22-Sep-2020 data refresh (01:29 -> 02:46) to the "Pais" sheet:
#  - numeric COVID counters refreshed for several countries
#  - "Santa Lucia"/"Timor Oriental" swap position (rows 204/205)
#  - "Montserrat"/"Islas Malvinas" swap position (rows 214/215)
#  - updated "datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 7045819
$ws.Range("C4").Value = 35975
$ws.Range("D4").Value = 4296840
$ws.Range("E4").Value = 2544507
$ws.Range("G4").Value = 354
$ws.Range("H4").Value = 204472

# --- Row 13: Argentina ---
$ws.Range("B13").Value = 640147
$ws.Range("C13").Value = 8782
$ws.Range("D13").Value = 508563
$ws.Range("E13").Value = 118102
$ws.Range("G13").Value = 429
$ws.Range("H13").Value = 13482

# --- Row 29: Canada ---
$ws.Range("B29").Value = 145415
$ws.Range("C29").Value = 1766
$ws.Range("D29").Value = 125534
$ws.Range("E29").Value = 10653

# --- Row 60: Chequia ---
$ws.Range("D60").Value = 25425
$ws.Range("E60").Value = 24817
$ws.Range("G60").Value = 19
$ws.Range("H60").Value = 522

# --- Row 72: Paraguay ---
$ws.Range("B72").Value = 34260
$ws.Range("C72").Value = 740
$ws.Range("D72").Value = 18629
$ws.Range("E72").Value = 14955
$ws.Range("G72").Value = 17
$ws.Range("H72").Value = 676

# --- Row 143: Mali ---
$ws.Range("B143").Value = 3024
$ws.Range("C143").Value = 11
$ws.Range("D143").Value = 2377
$ws.Range("E143").Value = 519

# --- Row 190: Monaco ---
$ws.Range("B190").Value = 195
$ws.Range("C190").Value = 2
$ws.Range("D190").Value = 159
$ws.Range("E190").Value = 35

# --- Rows 204/205: Santa Lucia and Timor Oriental swap places ---
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("B204").Value = 27
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 26
$ws.Range("E204").Value = 1
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

$ws.Range("A205").Value = "Timor Oriental"
$ws.Range("B205").Value = 27
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 27
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# --- Rows 214/215: Islas Malvinas and Montserrat swap places ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

# --- Updated timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 02:46"
